$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Edit the A1 title: remove the trailing newline from the Kyrgyz title text ---
$ws.Range("A1").Value = "8.3.1.2 Экономикадагы иш менен камсыз болгон бардык калктын чакан жана орто ишканаларда иштегендердин үлүшү"

# --- Add a new column N holding the 2023 data, matching the formatting of column M ---
$ws.Range("M3").Copy($ws.Range("N3"))

$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2023

$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 2.5449890821474286

$ws.Range("M6").Copy($ws.Range("N6"))
$ws.Range("N6").Value = 1.4569686017619159

# --- Row height tweaks (autofit-style adjustments that accompanied the edits) ---
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 17.25
$ws.Rows.Item(6).RowHeight = 17.25
